# Insert a new weekly price record for "Zanahoria" at row 276, shifting all
# subsequent rows down by one (old row 276..305 -> new row 277..306).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(276).Insert()

$ws.Cells.Item(276, 1).Value = 1
$ws.Cells.Item(276, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(276, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(276, 4).Value = 44769
$ws.Cells.Item(276, 5).Value = 15
$ws.Cells.Item(276, 6).Value = 100114013
$ws.Cells.Item(276, 7).Value = "Zanahoria"
$ws.Cells.Item(276, 8).Value = "Sin especificar"
$ws.Cells.Item(276, 9).Value = "Primera"
$ws.Cells.Item(276, 10).Value = 40
$ws.Cells.Item(276, 11).Value = 23000
$ws.Cells.Item(276, 12).Value = 24000
$ws.Cells.Item(276, 13).Value = 23500
$ws.Cells.Item(276, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(276, 15).Value = "Valle de Camiña"
$ws.Cells.Item(276, 16).Value = 940
$ws.Cells.Item(276, 17).Value = 25
$ws.Cells.Item(276, 18).Value = "Hortaliza"
